$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 17:50"

# Update country data rows (COVID-19 case counts refreshed; some rows swap
# which country occupies them because the table is ranked by total cases)
$ws.Range("A7").Value = "Alemania"
$ws.Range("B7").Value = 89451
$ws.Range("C7").Value = 4657
$ws.Range("D7").Value = 24575
$ws.Range("E7").Value = 63668
$ws.Range("F7").Value = 3936
$ws.Range("G7").Value = 101
$ws.Range("H7").Value = 1208

$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 11747
$ws.Range("C16").Value = 464
$ws.Range("D16").Value = 1979
$ws.Range("E16").Value = 9595
$ws.Range("F16").Value = 120
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 173

$ws.Range("A17").Value = "Austria"
$ws.Range("B17").Value = 11444
$ws.Range("C17").Value = 315
$ws.Range("D17").Value = 2022
$ws.Range("E17").Value = 9254
$ws.Range("F17").Value = 245
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 168

$ws.Range("A20").Value = "Brasil"
$ws.Range("B20").Value = 8195
$ws.Range("C20").Value = 151
$ws.Range("D20").Value = 127
$ws.Range("E20").Value = 7733
$ws.Range("F20").Value = 296
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 335

$ws.Range("A22").Value = "Suecia"
$ws.Range("B22").Value = 6131
$ws.Range("C22").Value = 563
$ws.Range("D22").Value = 205
$ws.Range("E22").Value = 5568
$ws.Range("F22").Value = 469
$ws.Range("G22").Value = 50
$ws.Range("H22").Value = 358

$ws.Range("A26").Value = "Chequia"
$ws.Range("B26").Value = 4091
$ws.Range("C26").Value = 233
$ws.Range("D26").Value = 71
$ws.Range("E26").Value = 3974
$ws.Range("F26").Value = 77
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 46

$ws.Range("A28").Value = "Dinamarca"
$ws.Range("B28").Value = 3757
$ws.Range("C28").Value = 371
$ws.Range("D28").Value = 1193
$ws.Range("E28").Value = 2425
$ws.Range("F28").Value = 153
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = 139

$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 3737
$ws.Range("C29").Value = 333
$ws.Range("D29").Value = 427
$ws.Range("E29").Value = 3288
$ws.Range("F29").Value = 31
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 22

$ws.Range("A31").Value = "Polonia"
$ws.Range("B31").Value = 3266
$ws.Range("C31").Value = 320
$ws.Range("D31").Value = 56
$ws.Range("E31").Value = 3145
$ws.Range("F31").Value = 50
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = 65

$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 3183
$ws.Range("C32").Value = 445
$ws.Range("D32").Value = 283
$ws.Range("E32").Value = 2778
$ws.Range("F32").Value = 83
$ws.Range("G32").Value = 7
$ws.Range("H32").Value = 122

$ws.Range("A33").Value = "Ecuador"
$ws.Range("B33").Value = 3163
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 65
$ws.Range("E33").Value = 2978
$ws.Range("F33").Value = 100
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 120

$ws.Range("A51").Value = "Argentina"
$ws.Range("B51").Value = 1265
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 266
$ws.Range("E51").Value = 960
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 39

$ws.Range("A72").Value = "Libano"
$ws.Range("B72").Value = 508
$ws.Range("C72").Value = 14
$ws.Range("D72").Value = 50
$ws.Range("E72").Value = 441
$ws.Range("F72").Value = 26
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 17

$ws.Range("A74").Value = "Tunez"
$ws.Range("B74").Value = 495
$ws.Range("C74").Value = 40
$ws.Range("D74").Value = 5
$ws.Range("E74").Value = 472
$ws.Range("F74").Value = 10
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 18

$ws.Range("A75").Value = "Letonia"
$ws.Range("B75").Value = 493
$ws.Range("C75").Value = 35
$ws.Range("D75").Value = 1
$ws.Range("E75").Value = 492
$ws.Range("F75").Value = 3
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0

$ws.Range("A76").Value = "Bulgaria"
$ws.Range("B76").Value = 485
$ws.Range("C76").Value = 28
$ws.Range("D76").Value = 30
$ws.Range("E76").Value = 441
$ws.Range("F76").Value = 18
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = 14

$ws.Range("A77").Value = "Kazajistan"
$ws.Range("B77").Value = 460
$ws.Range("C77").Value = 25
$ws.Range("D77").Value = 29
$ws.Range("E77").Value = 425
$ws.Range("F77").Value = 6
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 6

$ws.Range("A88").Value = "Reunion"
$ws.Range("B88").Value = 321
$ws.Range("C88").Value = 13
$ws.Range("D88").Value = 40
$ws.Range("E88").Value = 281
$ws.Range("F88").Value = 3
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0

$ws.Range("A90").Value = "Albania"
$ws.Range("B90").Value = 304
$ws.Range("C90").Value = 27
$ws.Range("D90").Value = 89
$ws.Range("E90").Value = 198
$ws.Range("F90").Value = 7
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 17

$ws.Range("A91").Value = "Burkina Faso"
$ws.Range("B91").Value = 302
$ws.Range("C91").Value = 14
$ws.Range("D91").Value = 50
$ws.Range("E91").Value = 236
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 16

$ws.Range("A92").Value = "Jordania"
$ws.Range("B92").Value = 299
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 45
$ws.Range("E92").Value = 249
$ws.Range("F92").Value = 5
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 5

$ws.Range("A94").Value = "Cuba"
$ws.Range("B94").Value = 269
$ws.Range("C94").Value = 36
$ws.Range("D94").Value = 15
$ws.Range("E94").Value = 248
$ws.Range("F94").Value = 7
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 6

$ws.Range("A95").Value = "Oman"
$ws.Range("B95").Value = 252
$ws.Range("C95").Value = 21
$ws.Range("D95").Value = 57
$ws.Range("E95").Value = 194
$ws.Range("F95").Value = 3
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 1

$ws.Range("A96").Value = "San Marino"
$ws.Range("B96").Value = 245
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 21
$ws.Range("E96").Value = 194
$ws.Range("F96").Value = 15
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 30

$ws.Range("A97").Value = "Vietnam"
$ws.Range("B97").Value = 237
$ws.Range("C97").Value = 4
$ws.Range("D97").Value = 85
$ws.Range("E97").Value = 152
$ws.Range("F97").Value = 3
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0

$ws.Range("A107").Value = "Estado de Palestina"
$ws.Range("B107").Value = 171
$ws.Range("C107").Value = 10
$ws.Range("D107").Value = 21
$ws.Range("E107").Value = 149
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 1

$ws.Range("A134").Value = "Guatemala"
$ws.Range("B134").Value = 50
$ws.Range("C134").Value = 3
$ws.Range("D134").Value = 12
$ws.Range("E134").Value = 37
$ws.Range("F134").Value = 1
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 1

$ws.Range("A135").Value = "Republica de Yibuti"
$ws.Range("B135").Value = 49
$ws.Range("C135").Value = 9
$ws.Range("D135").Value = 8
$ws.Range("E135").Value = 41
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0

$ws.Range("A136").Value = "Jamaica"
$ws.Range("B136").Value = 47
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 2
$ws.Range("E136").Value = 42
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 3

